# Commit: "Fruta / hortaliza, semanal"
# This edit inserts one new weekly price-report row into the daily logic
# sub-dataset sheet. The new row becomes row 26 (pushing the former rows
# 26..125 down to 27..126), and the sheet's used dimension grows from
# A1:T125 to A1:T126.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26; everything previously at/after row 26
# (through row 125) shifts down to rows 27..126, carrying its original
# values and formatting with it.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record.
$ws.Range("A26").Value = 8
$ws.Range("B26").Value = "Terminal La Palmera de La Serena"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 44715
$ws.Range("E26").Value = 4
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100109
$ws.Range("H26").Value = "Uva"
$ws.Range("I26").Value = 100109001
$ws.Range("J26").Value = "Uva"
$ws.Range("K26").Value = "Red Globe"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 400
$ws.Range("N26").Value = 9000
$ws.Range("O26").Value = 10000
$ws.Range("P26").Value = 9500
$ws.Range("Q26").Value = "$/bandeja 18 kilos"
$ws.Range("R26").Value = "Provincia del Elquí"
$ws.Range("S26").Value = 528
$ws.Range("T26").Value = 18
